$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 in place with the new record (concentration-check feature data).
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "aciclovir"
$ws.Range("C2").Value = "250mg/ml"
$ws.Range("D2").Value = "cimed"

# E2 must stay a text value ("143810181") rather than be coerced to a number.
# Build it in a scratch cell as a text formula, then paste-special as values
# into E2 so the result is a plain shared-string cell with no leftover
# number-format/style residue.
$scratch = $ws.Cells.Item(50, 50)
$scratch.Formula = '="143810181"'
$scratch.Copy()
$dest = $ws.Cells.Item(2, 5)
$dest.PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()
$excel.CutCopyMode = $false

$ws.Range("F2").Value = "Pendente"

# Rows 3-5 no longer apply; drop them so only the header + one record remain.
$ws.Range("A3:F5").Delete()
